# Update absPath and revisionPtr-driving metadata happen automatically on save;
# here we replicate the visible user actions: adding a third worksheet
# ("Planilha3") with the SQL-builder helper formulas, and leaving the
# selection/scroll state of Planilha2 the way the author left it before
# switching to the new sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Update Planilha2's view state (scrolled to row 18, G20:G29 selected) ---
$ws2.Activate()
$ws2.Range("G20:G29").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 18

# --- Add the new sheet right after Planilha2 ---
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "Planilha3"

# Column A width
$ws3.Columns.Item(1).ColumnWidth = 21.73

# Row 1: "Insert into produto (" label + full CONCAT formula
$ws3.Range("B1").Value = "Insert into produto ("
$ws3.Range("F1").Formula = "=CONCAT(B1:B18)"

# Rows 2-9: field name list + trailing comma helper formulas
$ws3.Range("A2").Value = "nome_Produto"
$ws3.Range("B2").Formula = '=A2&","'

$ws3.Range("A3").Value = "id_categoria_Produto"
$ws3.Range("A4").Value = "valorcusto_Produto"
$ws3.Range("A5").Value = "valorvenda_Produto"
$ws3.Range("A6").Value = "descricao_Produto"
$ws3.Range("A7").Value = "qtde_Produto"
$ws3.Range("A8").Value = "obs_Produto"
$ws3.Range("B3:B8").Formula = '=A3&","'

$ws3.Range("A9").Value = "status_Produto"
$ws3.Range("B9").Formula = "=A9"

# Row 10: ")values(" label
$ws3.Range("B10").Value = ")values("

# Rows 11-18: textbox/combobox names + quoted-value helper formulas
$ws3.Range("A11").Value = "txtNome"
$ws3.Range("B11").Formula = '="''""+"&A11&".Text+""'',"'

$ws3.Range("A12").Value = "cboIDCategoria"
$ws3.Range("A13").Value = "txtValorCusto"
$ws3.Range("A14").Value = "txtValorVenda"
$ws3.Range("A15").Value = "txtDescricao"
$ws3.Range("A16").Value = "txtQtde"
$ws3.Range("A17").Value = "txtOBS"
$ws3.Range("A18").Value = "cboStatus"
$ws3.Range("B12:B17").Formula = '="''""+"&A12&".Text+""'',"'

# B18 keeps the same pattern but with a closing parenthesis instead of a comma
$ws3.Range("B18").Formula = '="''""+"&A18&".Text+""'')"""'

# --- Final view state of the new sheet ---
$ws3.Range("F1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 190
